$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Actuarial Firm Summary")
$ws2 = $wb.Worksheets.Item("Actuarial Firm Summary 2")

# Merge "Segal Consulting" totals into "Segal" row (row 5) on sheet 1
$ws1.Range("B5").Value = 649806185770
$ws1.Range("C5").Value = 0.108115265762705
$ws1.Range("D5").Value = 160198525120

# Remove the now-obsolete "Segal Consulting" row (row 21); rows below shift up
$ws1.Rows(21).Delete()

# Apply the same Segal merge on sheet 2
$ws2.Range("B5").Value = 649806185770
$ws2.Range("C5").Value = 0.108115265762705
$ws2.Range("D5").Value = 160198525120

# Recompute the "Others" plug row on sheet 2 (row 17)
$ws2.Range("B17").Value = 107249819769
$ws2.Range("C17").Value = 0.017844309613008
$ws2.Range("D17").Value = 15056900803
